$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "wong3"

$ws.Range("B2").Value = 99
$ws.Range("C2").Value = 99
$ws.Range("E2").Value = 795
$ws.Range("F2").Value = 795
$ws.Range("B3").Value = 22
$ws.Range("C3").Value = 22
$ws.Range("E3").Value = 84
$ws.Range("F3").Value = 84
$ws.Range("B4").Value = 242
$ws.Range("C4").Value = 242
$ws.Range("E4").Value = 290
$ws.Range("F4").Value = 290
$ws.Range("B5").Value = 89
$ws.Range("C5").Value = 89
$ws.Range("E5").Value = 234
$ws.Range("F5").Value = 234
$ws.Range("H5").Value = 109
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 3
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 3
$ws.Range("E7").Value = 36
$ws.Range("F7").Value = 5
$ws.Range("B8").Value = 9
$ws.Range("C8").Value = 9
$ws.Range("E8").Value = 622
$ws.Range("F8").Value = 622
$ws.Range("B9").Value = 7
$ws.Range("C9").Value = 7
$ws.Range("E9").Value = 103
$ws.Range("F9").Value = 103
$ws.Range("B11").Value = 209
$ws.Range("C11").Value = 209
$ws.Range("E11").Value = 855
$ws.Range("F11").Value = 855
$ws.Range("B12").Value = 88
$ws.Range("C12").Value = 88
$ws.Range("E12").Value = 382
$ws.Range("F12").Value = 382
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 1
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("B14").Value = 67
$ws.Range("C14").Value = 67
$ws.Range("E14").Value = 271
$ws.Range("F14").Value = 271
$ws.Range("B15").Value = 43
$ws.Range("C15").Value = 43
$ws.Range("E15").Value = 115
$ws.Range("F15").Value = 112
$ws.Range("B16").Value = 281
$ws.Range("C16").Value = 281
$ws.Range("E16").Value = 580
$ws.Range("F16").Value = 580
$ws.Range("H16").Value = 19
$ws.Range("B17").Value = 214
$ws.Range("C17").Value = 214
$ws.Range("E17").Value = 857
$ws.Range("F17").Value = 857
$ws.Range("B19").Value = 38
$ws.Range("C19").Value = 38
$ws.Range("E19").Value = 367
$ws.Range("F19").Value = 367
$ws.Range("B20").Value = 2
$ws.Range("E20").Value = 38
$ws.Range("F20").Value = 5
$ws.Range("B22").Value = 4
$ws.Range("C22").Value = 4
$ws.Range("E22").Value = 4
$ws.Range("F22").Value = 4
$ws.Range("E23").Value = 109
$ws.Range("F23").Value = 109
$ws.Range("H23").Value = 78
$ws.Range("B24").Value = 3
$ws.Range("C24").Value = 3
$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 3
$ws.Range("B25").Value = 329
$ws.Range("C25").Value = 329
$ws.Range("E25").Value = 631
$ws.Range("F25").Value = 631
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 1
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("B27").Value = 85
$ws.Range("C27").Value = 85
$ws.Range("E27").Value = 471
$ws.Range("F27").Value = 471
$ws.Range("B29").Value = 4
$ws.Range("C29").Value = 4
$ws.Range("E29").Value = 32
$ws.Range("F29").Value = 32
$ws.Range("H29").Value = 46
$ws.Range("E31").Value = 231
$ws.Range("F31").Value = 231
$ws.Range("H31").Value = 15
$ws.Range("B32").Value = 18
$ws.Range("C32").Value = 18
$ws.Range("E32").Value = 679
$ws.Range("F32").Value = 679
$ws.Range("B33").Value = 66
$ws.Range("C33").Value = 66
$ws.Range("E33").Value = 240
$ws.Range("F33").Value = 240
$ws.Range("B34").Value = 2
$ws.Range("C34").Value = 2
$ws.Range("E34").Value = 3
$ws.Range("F34").Value = 3
$ws.Range("B35").Value = 2
$ws.Range("C35").Value = 2
$ws.Range("E35").Value = 71
$ws.Range("F35").Value = 71
$ws.Range("B36").Value = 32
$ws.Range("C36").Value = 32
$ws.Range("E36").Value = 856
$ws.Range("F36").Value = 856
$ws.Range("B39").Value = 13
$ws.Range("C39").Value = 13
$ws.Range("E39").Value = 242
$ws.Range("F39").Value = 242
$ws.Range("B41").Value = 3
$ws.Range("E41").Value = 32
$ws.Range("B43").Value = 7
$ws.Range("C43").Value = 7
$ws.Range("E43").Value = 32
$ws.Range("F43").Value = 32
$ws.Range("B45").Value = 170
$ws.Range("C45").Value = 170
$ws.Range("E45").Value = 326
$ws.Range("F45").Value = 326
$ws.Range("B46").Value = 48
$ws.Range("C46").Value = 48
$ws.Range("E46").Value = 775
$ws.Range("F46").Value = 775
$ws.Range("B47").Value = 27
$ws.Range("C47").Value = 27
$ws.Range("E47").Value = 88
$ws.Range("F47").Value = 88
$ws.Range("B48").Value = 1
$ws.Range("C48").Value = 1
$ws.Range("E48").Value = 1
$ws.Range("F48").Value = 1
$ws.Range("B49").Value = 10
$ws.Range("C49").Value = 10
$ws.Range("E49").Value = 131
$ws.Range("F49").Value = 131
$ws.Range("B50").Value = 16
$ws.Range("C50").Value = 16
$ws.Range("E50").Value = 591
$ws.Range("F50").Value = 591
$ws.Range("B51").Value = 8
$ws.Range("C51").Value = 1
$ws.Range("E51").Value = 183
$ws.Range("F51").Value = 182
$ws.Range("H51").Value = 21
$ws.Range("B52").Value = 9
$ws.Range("C52").Value = 9
$ws.Range("E52").Value = 516
$ws.Range("F52").Value = 513
$ws.Range("E53").Value = 49
$ws.Range("F53").Value = 47
$ws.Range("H53").Value = 19
$ws.Range("B55").Value = 33
$ws.Range("C55").Value = 33
$ws.Range("E55").Value = 858
$ws.Range("F55").Value = 858
$ws.Range("B56").Value = 4
$ws.Range("C56").Value = 4
$ws.Range("E56").Value = 4
$ws.Range("F56").Value = 4
$ws.Range("B57").Value = 86
$ws.Range("C57").Value = 86
$ws.Range("E57").Value = 564
$ws.Range("F57").Value = 564
$ws.Range("B58").Value = 18
$ws.Range("C58").Value = 18
$ws.Range("E58").Value = 423
$ws.Range("F58").Value = 423
$ws.Range("B59").Value = 2
$ws.Range("C59").Value = 2
$ws.Range("E59").Value = 1
$ws.Range("F59").Value = 1
$ws.Range("B60").Value = 11
$ws.Range("C60").Value = 11
$ws.Range("E60").Value = 262
$ws.Range("F60").Value = 262
$ws.Range("B62").Value = 97
$ws.Range("C62").Value = 97
$ws.Range("E62").Value = 845
$ws.Range("F62").Value = 845
$ws.Range("B63").Value = 28
$ws.Range("C63").Value = 28
$ws.Range("E63").Value = 852
$ws.Range("F63").Value = 852
$ws.Range("B64").Value = 86
$ws.Range("C64").Value = 86
$ws.Range("E64").Value = 617
$ws.Range("F64").Value = 617
$ws.Range("B66").Value = 87
$ws.Range("C66").Value = 87
$ws.Range("E66").Value = 661
$ws.Range("F66").Value = 660
$ws.Range("B67").Value = 2
$ws.Range("C67").Value = 2
$ws.Range("E67").Value = 3
$ws.Range("F67").Value = 3
$ws.Range("B68").Value = 286
$ws.Range("C68").Value = 286
$ws.Range("E68").Value = 254
$ws.Range("F68").Value = 254
$ws.Range("B69").Value = 1
$ws.Range("C69").Value = 1
$ws.Range("E69").Value = 1
$ws.Range("F69").Value = 1
$ws.Range("B70").Value = 25
$ws.Range("C70").Value = 25
$ws.Range("E70").Value = 690
$ws.Range("F70").Value = 690
$ws.Range("B74").Value = 213
$ws.Range("C74").Value = 213
$ws.Range("E74").Value = 928
$ws.Range("F74").Value = 928
$ws.Range("B75").Value = 15
$ws.Range("C75").Value = 15
$ws.Range("E75").Value = 177
$ws.Range("F75").Value = 177
$ws.Range("B76").Value = 3
$ws.Range("C76").Value = 3
$ws.Range("E76").Value = 27
$ws.Range("F76").Value = 27
$ws.Range("B79").Value = 103
$ws.Range("C79").Value = 101
$ws.Range("E79").Value = 124
$ws.Range("F79").Value = 121
$ws.Range("B80").Value = 116
$ws.Range("C80").Value = 116
$ws.Range("E80").Value = 168
$ws.Range("F80").Value = 168
$ws.Range("B81").Value = 144
$ws.Range("C81").Value = 144
$ws.Range("E81").Value = 372
$ws.Range("F81").Value = 372
$ws.Range("B82").Value = 40
$ws.Range("C82").Value = 40
$ws.Range("E82").Value = 568
$ws.Range("F82").Value = 566
$ws.Range("B83").Value = 3
$ws.Range("E83").Value = 33
$ws.Range("B84").Value = 56
$ws.Range("C84").Value = 56
$ws.Range("E84").Value = 650
$ws.Range("F84").Value = 650
$ws.Range("B85").Value = 14
$ws.Range("C85").Value = 14
$ws.Range("E85").Value = 192
$ws.Range("F85").Value = 192
$ws.Range("B86").Value = 217
$ws.Range("C86").Value = 217
$ws.Range("E86").Value = 861
$ws.Range("F86").Value = 861
$ws.Range("B89").Value = 23
$ws.Range("C89").Value = 23
$ws.Range("E89").Value = 373
$ws.Range("F89").Value = 373
$ws.Range("B90").Value = 1
$ws.Range("C90").Value = 1
$ws.Range("B91").Value = 103
$ws.Range("C91").Value = 103
$ws.Range("E91").Value = 234
$ws.Range("F91").Value = 234
$ws.Range("H91").Value = 49
$ws.Range("B92").Value = 32
$ws.Range("C92").Value = 32
$ws.Range("E92").Value = 855
$ws.Range("F92").Value = 855
$ws.Range("B93").Value = 167
$ws.Range("C93").Value = 167
$ws.Range("E93").Value = 248
$ws.Range("F93").Value = 248
